$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45189 -> 45190, i.e. 2023-09-20 -> 2023-09-21) for every data row (2..509).
$row = 2
while ($true) {
    $cell = $ws.Cells.Item($row, 3)
    $val = $cell.Value2
    if ($val -eq $null) {
        break
    }
    if ($val -eq 45189) {
        $cell.Value2 = 45190
    }
    $row = $row + 1
}
